# Update the "Förändrad" date column (C) for all data rows.
# Every value in C2:C387 changes from serial date 45204 (2023-10-05)
# to serial date 45205 (2023-10-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 387 }

$range = $ws.Range("C2:C$lastRow")
foreach ($cell in $range.Cells) {
    if ($cell.Value2 -eq 45204) {
        $cell.Value = 45205
    }
}
